$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 192.83333
$ws.Range("J42").Value = 265
$ws.Range("L42").Value = 795
$ws.Range("N42").Value = -1255

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H88").Value = 991.4666999999999
$ws.Range("J88").Value = 1137.9
$ws.Range("L88").Value = 1137.9
$ws.Range("N88").Value = -1949.9

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H91").Value = 991.4666999999999
$ws.Range("J91").Value = 1137.9
$ws.Range("L91").Value = 1137.9
$ws.Range("N91").Value = -3945.9

$ws.Range("H132").Value = 2637.6035
$ws.Range("I132").Value = 2461.2263
$ws.Range("K132").Value = 7383.678899999999
$ws.Range("M132").Value = -4853.678899999999

$ws.Range("H135").Value = 6676.913
$ws.Range("I135").Value = 8345.375
$ws.Range("J135").Value = 2863.2856
$ws.Range("K135").Value = 75108.375
$ws.Range("L135").Value = 25769.5704
$ws.Range("M135").Value = -72573.375
$ws.Range("N135").Value = -30839.5704

$ws.Range("H138").Value = 2820.9387
$ws.Range("J138").Value = 3214.0833
$ws.Range("L138").Value = 9642.249899999999
$ws.Range("N138").Value = -19922.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6491.8833
$ws.Range("I32").Value = 6214.527
$ws.Range("K32").Value = 6214.527
$ws.Range("M32").Value = -5927.527

$ws.Range("H61").Value = 11281.046
$ws.Range("J61").Value = 6999.75
$ws.Range("L61").Value = 6999.75
$ws.Range("N61").Value = -7423.75

$ws.Range("H74").Value = 8706.267
$ws.Range("I74").Value = 17415.5
$ws.Range("K74").Value = 17415.5
$ws.Range("M74").Value = -16541.5

$ws.Range("H77").Value = 8706.267
$ws.Range("I77").Value = 17415.5
$ws.Range("K77").Value = 87077.5
$ws.Range("M77").Value = -82709.5

$ws.Range("H102").Value = 17808.459
$ws.Range("I102").Value = 27875.5
$ws.Range("J102").Value = 7741.4165
$ws.Range("K102").Value = 27875.5
$ws.Range("L102").Value = 7741.4165
$ws.Range("M102").Value = -26253.5
$ws.Range("N102").Value = -10985.4165

$ws.Range("H132").Value = 2879.7273
$ws.Range("I132").Value = 1975.4193
$ws.Range("J132").Value = 5036.154
$ws.Range("K132").Value = 5926.257900000001
$ws.Range("L132").Value = 15108.462
$ws.Range("M132").Value = -3396.257900000001
$ws.Range("N132").Value = -20168.462

$ws.Range("H136").Value = 11281.046
$ws.Range("J136").Value = 6999.75
$ws.Range("L136").Value = 20999.25
$ws.Range("N136").Value = -26099.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 15194.619
$ws.Range("I99").Value = 19442.285
$ws.Range("J99").Value = 6699.2856
$ws.Range("K99").Value = 19442.285
$ws.Range("L99").Value = 6699.2856
$ws.Range("M99").Value = -17944.285
$ws.Range("N99").Value = -9695.285599999999

$ws.Range("H134").Value = 8255.308000000001
$ws.Range("I134").Value = 8538.362999999999
$ws.Range("K134").Value = 25615.089
$ws.Range("M134").Value = -23080.089

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 59877.53
$ws.Range("I16").Value = 1060.4
$ws.Range("J16").Value = 143902
$ws.Range("K16").Value = 1060.4
$ws.Range("L16").Value = 143902
$ws.Range("M16").Value = -773.4000000000001
$ws.Range("N16").Value = -144476

$ws.Range("H31").Value = 7212.425
$ws.Range("I31").Value = 8244.6
$ws.Range("K31").Value = 8244.6
$ws.Range("M31").Value = -7949.6

$ws.Range("H34").Value = 7212.425
$ws.Range("I34").Value = 8244.6
$ws.Range("K34").Value = 8244.6
$ws.Range("M34").Value = -8042.6

$ws.Range("H105").Value = 124637.06
$ws.Range("I105").Value = 162596.31
$ws.Range("K105").Value = 162596.31
$ws.Range("M105").Value = -160849.31

$ws.Range("H113").Value = 59877.53
$ws.Range("I113").Value = 1060.4
$ws.Range("J113").Value = 143902
$ws.Range("K113").Value = 1060.4
$ws.Range("L113").Value = 143902
$ws.Range("M113").Value = 1109.6
$ws.Range("N113").Value = -148242

$ws.Range("H122").Value = 12188
$ws.Range("I122").Value = 15989.875
$ws.Range("K122").Value = 47969.625
$ws.Range("M122").Value = -45519.625

$ws.Range("H132").Value = 1835.7
$ws.Range("I132").Value = 1588.4
$ws.Range("K132").Value = 4765.200000000001
$ws.Range("M132").Value = -2235.200000000001

$ws.Range("H141").Value = 241805.4
$ws.Range("J141").Value = 255149.84
$ws.Range("L141").Value = 255149.84
$ws.Range("N141").Value = -265509.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 51338596
$ws.Range("I4").Value = 41080860
$ws.Range("K4").Value = 123242580
$ws.Range("M4").Value = -123242468

$ws.Range("H11").Value = 77197.62
$ws.Range("I11").Value = 147.8
$ws.Range("J11").Value = 125353.75
$ws.Range("K11").Value = 443.4
$ws.Range("L11").Value = 376061.25
$ws.Range("M11").Value = -303.4
$ws.Range("N11").Value = -376341.25

$ws.Range("H80").Value = 156933
$ws.Range("J80").Value = 87986.3
$ws.Range("L80").Value = 263958.9
$ws.Range("N80").Value = -265830.9

$ws.Range("H83").Value = 156933
$ws.Range("J83").Value = 87986.3
$ws.Range("L83").Value = 791876.7000000001
$ws.Range("N83").Value = -801236.7000000001

$ws.Range("H131").Value = 11113830
$ws.Range("J131").Value = 2023.2152
$ws.Range("L131").Value = 6069.6456
$ws.Range("N131").Value = -16149.6456

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7447.1904
$ws.Range("I102").Value = 8946.375
$ws.Range("K102").Value = 8946.375
$ws.Range("M102").Value = -7324.375

$ws.Range("H113").Value = 9895.214
$ws.Range("I113").Value = 12453.4
$ws.Range("J113").Value = 3499.75
$ws.Range("K113").Value = 12453.4
$ws.Range("L113").Value = 3499.75
$ws.Range("M113").Value = -10283.4
$ws.Range("N113").Value = -7839.75

$ws.Range("H123").Value = 12332.6
$ws.Range("J123").Value = 12332.6
$ws.Range("L123").Value = 12332.6
$ws.Range("N123").Value = -17232.6

$ws.Range("H132").Value = 4600.095
$ws.Range("I132").Value = 4630.1
$ws.Range("K132").Value = 13890.3
$ws.Range("M132").Value = -11360.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1389.5333
$ws.Range("I55").Value = 460.3
$ws.Range("J55").Value = 3248
$ws.Range("K55").Value = 460.3
$ws.Range("L55").Value = 3248
$ws.Range("M55").Value = -287.3
$ws.Range("N55").Value = -3594

$ws.Range("H122").Value = 6495.7407
$ws.Range("I122").Value = 6112.533
$ws.Range("K122").Value = 18337.599
$ws.Range("M122").Value = -15887.599

$ws.Range("H132").Value = 598955.25
$ws.Range("I132").Value = 1065709.1
$ws.Range("J132").Value = 4904.8184
$ws.Range("K132").Value = 3197127.3
$ws.Range("L132").Value = 14714.4552
$ws.Range("M132").Value = -3194597.3
$ws.Range("N132").Value = -19774.4552

$ws.Range("H136").Value = 6305.1904
$ws.Range("I136").Value = 5799.5
$ws.Range("J136").Value = 6424.1763
$ws.Range("K136").Value = 17398.5
$ws.Range("L136").Value = 19272.5289
$ws.Range("M136").Value = -14848.5
$ws.Range("N136").Value = -24372.5289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 80069.836
$ws.Range("I46").Value = 73994
$ws.Range("J46").Value = 81285
$ws.Range("K46").Value = 73994
$ws.Range("L46").Value = 81285
$ws.Range("M46").Value = -73763
$ws.Range("N46").Value = -81747

$ws.Range("H107").Value = 45642.57
$ws.Range("I107").Value = 4499.5
$ws.Range("J107").Value = 100500
$ws.Range("K107").Value = 13498.5
$ws.Range("L107").Value = 301500
$ws.Range("M107").Value = -11578.5
$ws.Range("N107").Value = -305340

$ws.Range("H113").Value = 1009.4545
$ws.Range("I113").Value = 548.9677
$ws.Range("K113").Value = 1646.9031
$ws.Range("M113").Value = 523.0969

$ws.Range("H122").Value = 18866.545
$ws.Range("I122").Value = 1946.579
$ws.Range("J122").Value = 41829.355
$ws.Range("K122").Value = 5839.737
$ws.Range("L122").Value = 125488.065
$ws.Range("M122").Value = -3389.737
$ws.Range("N122").Value = -130388.065

$ws.Range("H132").Value = 7841.145
$ws.Range("J132").Value = 5081.2856
$ws.Range("L132").Value = 15243.8568
$ws.Range("N132").Value = -20303.8568

$ws.Range("H134").Value = 80069.836
$ws.Range("I134").Value = 73994
$ws.Range("J134").Value = 81285
$ws.Range("K134").Value = 221982
$ws.Range("L134").Value = 243855
$ws.Range("M134").Value = -219447
$ws.Range("N134").Value = -248925
